{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the \"Features \" paragraph so we can insert the new paragraph right after it.\nlet featuresParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Features\") {\n    featuresParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!featuresParagraph) {\n  throw new Error(\"Could not find 'Features' paragraph\");\n}\n\nconst newParagraph = featuresParagraph.insertParagraph(\"Following are the features\", Word.InsertLocation.after);\nnewParagraph.alignment = Word.Alignment.justified;\nnewParagraph.font.bold = true;\nnewParagraph.font.size = 12;\n\nawait context.sync();\n", "ps1": "$doc = $word.ActiveDocument\n\n# Locate the \"Features\" paragraph.\n$searchRange = $doc.Content\n$null = $searchRange.Find.Execute(\"Features\")\n$featuresParagraph = $searchRange.Paragraphs(1)\n\n# Insert a new paragraph right after it (inherits the \"Features\" paragraph's\n# formatting: bold, size 24/24 half-points, justified).\n$null = $featuresParagraph.Range.InsertParagraphAfter()\n$newParagraph = $featuresParagraph.Next()\n$newParagraph.Range.Text = \"Following are the features\"\n"}
